# Update the "想去人数" (interested-count) figures in column F for both the
# "展览" sheet (rows 3-14) and the "全部类型" sheet (rows 4-15, one row lower
# because it has an extra data row before the matching entries).

$wb = $excel.ActiveWorkbook

$newValues = @(48, 1437, 331, 1046, 10808, 25, 85, 298, 1048, 726, 12132, 12585)

# Sheet "展览" -> rows 3..14
$wsExpo = $wb.Worksheets.Item("展览")
for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = 3 + $i
    $wsExpo.Cells.Item($row, 6).Value = $newValues[$i]
}

# Sheet "全部类型" -> rows 4..15
$wsAll = $wb.Worksheets.Item("全部类型")
for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = 4 + $i
    $wsAll.Cells.Item($row, 6).Value = $newValues[$i]
}
